$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 607 ("2026/12/29" block),
# shifting all rows from 607 downward by 2 (607->609 ... 648->650).
$ws.Rows("607:608").Insert()

# The new rows hold two additional hourly readings for 2026/01/10 (Sat),
# continuing the existing sequence that already fills rows 602-606
# (02:00, 06:00, 10:00, 12:00, 15:00) with 18:00 and 20:00.

# Column A holds the date as plain text (e.g. "2026/01/10"), not a real
# date value, so format the cells as Text first to stop Excel from
# auto-converting the string into a date serial, then restore the
# default "Normal" style so no stray custom number format is left
# behind in the saved workbook.
$ws.Range("A607:A608").NumberFormat = "@"

$ws.Range("A607").Value = "2026/01/10"
$ws.Range("B607").Value = "土"
$ws.Range("C607").Value = 18
$ws.Range("D607").Value = 201

$ws.Range("A608").Value = "2026/01/10"
$ws.Range("B608").Value = "土"
$ws.Range("C608").Value = 20
$ws.Range("D608").Value = 201

$ws.Range("A607:A608").Style = "Normal"
